# Update cryptocurrency price list (Price and Volume(1h) columns),
# including 3 pairs of rows whose coin data (rank order) swapped positions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.116.96'
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.986.35'
$ws.Range("E3").Value = '  -3.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '536.44'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.53'
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '2.980.48'
$ws.Range("E8").Value = '  -3.04%  '
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.148'
$ws.Range("E10").Value = '  -4.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.08'
$ws.Range("E11").Value = '  -0.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.446'
$ws.Range("E12").Value = '  -2.38%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000221'
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.93'
$ws.Range("E14").Value = '  -1.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.461.55'
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.052.61'
$ws.Range("E17").Value = '  -3.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.985.52'
$ws.Range("E18").Value = '  -2.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.59'
$ws.Range("E19").Value = '  -1.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '465.36'
$ws.Range("E20").Value = '  -4.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.20'
$ws.Range("E21").Value = '  -1.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.672'
$ws.Range("E22").Value = '  -3.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.92'
$ws.Range("E23").Value = '  -3.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.51'
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.96'
$ws.Range("E25").Value = '  -2.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.995'
$ws.Range("E26").Value = '  -0.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.67'
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("E28").Value = '  -5.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  +3.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.88'
$ws.Range("E31").Value = '  -0.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '25.48'
$ws.Range("E32").Value = '  -2.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '55.32'
$ws.Range("B34").Value = 'Stacks'
$ws.Range("C34").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.28'
$ws.Range("E34").Value = '  -4.89%  '
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.41'
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.85'
$ws.Range("E36").Value = '  -3.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '452.56'
$ws.Range("E37").Value = '  -7.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.173.51'
$ws.Range("E38").Value = '  +0.36%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0785'
$ws.Range("E39").Value = '  -1.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0383'
$ws.Range("E40").Value = '  -3.59%  '
$ws.Range("E41").Value = '  +0.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.10'
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.47'
$ws.Range("E43").Value = '  -6.84%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.31'
$ws.Range("E44").Value = '  +10.94%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.243'
$ws.Range("E46").Value = '  -4.34%  '
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.00'
$ws.Range("E47").Value = '  -2.00%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '119.28'
$ws.Range("E48").Value = '  -1.99%  '
$ws.Range("E49").Value = '  -0.96%  '
$ws.Range("E50").Value = '  -8.60%  '
$ws.Range("E51").Value = '  +6.63%  '
